$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 28, shifting old rows 28 and 29 down to 29 and 30
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new weekly entry
$ws.Range("A28").Value = 11
$ws.Range("B28").Value = "Vega Monumental Concepción"
$ws.Range("C28").Value = "Bíobío"
$ws.Range("D28").Value = 44769
$ws.Range("D28").NumberFormat = $ws.Range("D29").NumberFormat
$ws.Range("E28").Value = 8
$ws.Range("F28").Value = 100112026
$ws.Range("G28").Value = "Haba"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 100
$ws.Range("K28").Value = 18000
$ws.Range("L28").Value = 20000
$ws.Range("M28").Value = 19000
$ws.Range("N28").Value = "`$/saco 25 kilos"
$ws.Range("O28").Value = "Región de Coquimbo"
$ws.Range("P28").Value = 760
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = "Hortaliza"
